$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3361.8108
$ws.Range("I64").Value = 3497.6667
$ws.Range("J64").Value = 2995
$ws.Range("K64").Value = 3497.6667
$ws.Range("L64").Value = 2995
$ws.Range("M64").Value = -3249.6667
$ws.Range("N64").Value = -3491
$ws.Range("H67").Value = 3361.8108
$ws.Range("I67").Value = 3497.6667
$ws.Range("J67").Value = 2995
$ws.Range("K67").Value = 3497.6667
$ws.Range("L67").Value = 2995
$ws.Range("M67").Value = -2639.6667
$ws.Range("N67").Value = -4711
$ws.Range("H96").Value = 308.4
$ws.Range("I96").Value = 289.7143
$ws.Range("J96").Value = 352
$ws.Range("K96").Value = 869.1428999999999
$ws.Range("L96").Value = 1056
$ws.Range("M96").Value = 503.8571000000001
$ws.Range("N96").Value = -3802
$ws.Range("H125").Value = 4543.75
$ws.Range("I125").Value = 8188
$ws.Range("J125").Value = 3702.7693
$ws.Range("K125").Value = 73692
$ws.Range("L125").Value = 33324.9237
$ws.Range("M125").Value = -71232
$ws.Range("N125").Value = -38244.9237
$ws.Range("H137").Value = 1108.7671
$ws.Range("I137").Value = 904.918
$ws.Range("J137").Value = 2145
$ws.Range("K137").Value = 2714.754
$ws.Range("L137").Value = 6435
$ws.Range("M137").Value = -164.7539999999999
$ws.Range("N137").Value = -11535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3092.22
$ws.Range("I32").Value = 2037.2727
$ws.Range("J32").Value = 10828.5
$ws.Range("K32").Value = 2037.2727
$ws.Range("L32").Value = 10828.5
$ws.Range("M32").Value = -1750.2727
$ws.Range("N32").Value = -11402.5
$ws.Range("H45").Value = 9369.583000000001
$ws.Range("I45").Value = 15333.571
$ws.Range("J45").Value = 1020
$ws.Range("K45").Value = 15333.571
$ws.Range("L45").Value = 1020
$ws.Range("M45").Value = -14956.571
$ws.Range("N45").Value = -1774
$ws.Range("H61").Value = 3699.4048
$ws.Range("I61").Value = 3963.75
$ws.Range("J61").Value = 2113.3333
$ws.Range("K61").Value = 3963.75
$ws.Range("L61").Value = 2113.3333
$ws.Range("M61").Value = -3751.75
$ws.Range("N61").Value = -2537.3333
$ws.Range("H74").Value = 1336.7667
$ws.Range("I74").Value = 1401.1177
$ws.Range("J74").Value = 1252.6154
$ws.Range("K74").Value = 1401.1177
$ws.Range("L74").Value = 1252.6154
$ws.Range("M74").Value = -527.1177
$ws.Range("N74").Value = -3000.6154
$ws.Range("H77").Value = 1336.7667
$ws.Range("I77").Value = 1401.1177
$ws.Range("J77").Value = 1252.6154
$ws.Range("K77").Value = 7005.5885
$ws.Range("L77").Value = 6263.076999999999
$ws.Range("M77").Value = -2637.5885
$ws.Range("N77").Value = -14999.077
$ws.Range("H102").Value = 12347279
$ws.Range("I102").Value = 12347279
$ws.Range("K102").Value = 12347279
$ws.Range("M102").Value = -12345657
$ws.Range("H110").Value = 2925
$ws.Range("J110").Value = 3500
$ws.Range("L110").Value = 3500
$ws.Range("N110").Value = -7590
$ws.Range("H122").Value = 4279240
$ws.Range("I122").Value = 4279240
$ws.Range("K122").Value = 12837720
$ws.Range("M122").Value = -12835270
$ws.Range("H132").Value = 2567089.2
$ws.Range("I132").Value = 2326.84
$ws.Range("J132").Value = 7147022
$ws.Range("K132").Value = 6980.52
$ws.Range("L132").Value = 21441066
$ws.Range("M132").Value = -4450.52
$ws.Range("N132").Value = -21446126
$ws.Range("H136").Value = 3699.4048
$ws.Range("I136").Value = 3963.75
$ws.Range("J136").Value = 2113.3333
$ws.Range("K136").Value = 11891.25
$ws.Range("L136").Value = 6339.999899999999
$ws.Range("M136").Value = -9341.25
$ws.Range("N136").Value = -11439.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19373.5
$ws.Range("I20").Value = 1488
$ws.Range("J20").Value = 37259
$ws.Range("K20").Value = 1488
$ws.Range("L20").Value = 37259
$ws.Range("M20").Value = -1241
$ws.Range("N20").Value = -37753
$ws.Range("H86").Value = 1956.1111
$ws.Range("I86").Value = 1934.1666
$ws.Range("K86").Value = 1934.1666
$ws.Range("M86").Value = -811.1666
$ws.Range("H89").Value = 1956.1111
$ws.Range("I89").Value = 1934.1666
$ws.Range("K89").Value = 9670.833000000001
$ws.Range("M89").Value = -4054.833000000001
$ws.Range("H94").Value = 2150.9
$ws.Range("I94").Value = 754.5
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 754.5
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -303.5
$ws.Range("N94").Value = -3402
$ws.Range("H103").Value = 43000
$ws.Range("J103").Value = 43000
$ws.Range("L103").Value = 43000
$ws.Range("N103").Value = -45344
$ws.Range("H105").Value = 29425418
$ws.Range("I105").Value = 50021530
$ws.Range("J105").Value = 2401.5715
$ws.Range("K105").Value = 50021530
$ws.Range("L105").Value = 2401.5715
$ws.Range("M105").Value = -50019783
$ws.Range("N105").Value = -5895.5715
$ws.Range("H107").Value = 142858610
$ws.Range("I107").Value = 166668050
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 166668050
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -166666130
$ws.Range("N107").Value = -5840
$ws.Range("H134").Value = 2919.8708
$ws.Range("I134").Value = 2966.0393
$ws.Range("J134").Value = 2705.818
$ws.Range("K134").Value = 8898.117899999999
$ws.Range("L134").Value = 8117.454000000001
$ws.Range("M134").Value = -6363.117899999999
$ws.Range("N134").Value = -13187.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6397.537
$ws.Range("I31").Value = 1739.093
$ws.Range("J31").Value = 14743.917
$ws.Range("K31").Value = 1739.093
$ws.Range("L31").Value = 14743.917
$ws.Range("M31").Value = -1444.093
$ws.Range("N31").Value = -15333.917
$ws.Range("H34").Value = 6397.537
$ws.Range("I34").Value = 1739.093
$ws.Range("J34").Value = 14743.917
$ws.Range("K34").Value = 1739.093
$ws.Range("L34").Value = 14743.917
$ws.Range("M34").Value = -1537.093
$ws.Range("N34").Value = -15147.917
$ws.Range("H58").Value = 931.6479
$ws.Range("I58").Value = 564.13464
$ws.Range("K58").Value = 564.13464
$ws.Range("M58").Value = -361.13464
$ws.Range("H132").Value = 2711.3914
$ws.Range("I132").Value = 2066.1538
$ws.Range("J132").Value = 3550.2
$ws.Range("K132").Value = 6198.4614
$ws.Range("L132").Value = 10650.6
$ws.Range("M132").Value = -3668.4614
$ws.Range("N132").Value = -15710.6
$ws.Range("H136").Value = 931.6479
$ws.Range("I136").Value = 564.13464
$ws.Range("K136").Value = 1692.40392
$ws.Range("M136").Value = 857.59608
$ws.Range("H141").Value = 33765.332
$ws.Range("J141").Value = 33765.332
$ws.Range("L141").Value = 33765.332
$ws.Range("N141").Value = -44125.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 366
$ws.Range("I92").Value = 274.33334
$ws.Range("J92").Value = 421
$ws.Range("K92").Value = 823.0000200000001
$ws.Range("L92").Value = 1263
$ws.Range("M92").Value = 424.9999799999999
$ws.Range("N92").Value = -3759
$ws.Range("H117").Value = 20843832
$ws.Range("J117").Value = 33339908
$ws.Range("L117").Value = 100019724
$ws.Range("N117").Value = -100026608
$ws.Range("H121").Value = 975.7931
$ws.Range("I121").Value = 482.5
$ws.Range("J121").Value = 1054.72
$ws.Range("K121").Value = 1447.5
$ws.Range("L121").Value = 3164.16
$ws.Range("M121").Value = -137.5
$ws.Range("N121").Value = -5784.16
$ws.Range("H129").Value = 1196.6
$ws.Range("I129").Value = 720
$ws.Range("J129").Value = 1741.2858
$ws.Range("K129").Value = 2160
$ws.Range("L129").Value = 5223.857400000001
$ws.Range("M129").Value = 2840
$ws.Range("N129").Value = -15223.8574
$ws.Range("H136").Value = 12052.5
$ws.Range("I136").Value = 50315
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 150945
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -145845
$ws.Range("N136").Value = -23400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6052.815
$ws.Range("I126").Value = 11671.1
$ws.Range("K126").Value = 35013.3
$ws.Range("M126").Value = -32543.3
$ws.Range("H132").Value = 2475.3704
$ws.Range("I132").Value = 1797.0667
$ws.Range("J132").Value = 3323.25
$ws.Range("K132").Value = 5391.2001
$ws.Range("L132").Value = 9969.75
$ws.Range("M132").Value = -2861.2001
$ws.Range("N132").Value = -15029.75
$ws.Range("H139").Value = 67548.78
$ws.Range("J139").Value = 67548.78
$ws.Range("L139").Value = 67548.78
$ws.Range("N139").Value = -77828.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1912.591
$ws.Range("I7").Value = 1756.9286
$ws.Range("J7").Value = 2185
$ws.Range("K7").Value = 1756.9286
$ws.Range("L7").Value = 2185
$ws.Range("M7").Value = -1644.9286
$ws.Range("N7").Value = -2409
$ws.Range("H82").Value = 291658.72
$ws.Range("I82").Value = 455940.6
$ws.Range("J82").Value = 65771.125
$ws.Range("K82").Value = 455940.6
$ws.Range("L82").Value = 65771.125
$ws.Range("M82").Value = -455579.6
$ws.Range("N82").Value = -66493.125
$ws.Range("H85").Value = 291658.72
$ws.Range("I85").Value = 455940.6
$ws.Range("J85").Value = 65771.125
$ws.Range("K85").Value = 455940.6
$ws.Range("L85").Value = 65771.125
$ws.Range("M85").Value = -454692.6
$ws.Range("N85").Value = -68267.125
$ws.Range("H93").Value = 71429280
$ws.Range("I93").Value = 666.6667
$ws.Range("J93").Value = 125000744
$ws.Range("K93").Value = 666.6667
$ws.Range("L93").Value = 125000744
$ws.Range("M93").Value = 581.3333
$ws.Range("N93").Value = -125003240
$ws.Range("H126").Value = 1912.591
$ws.Range("I126").Value = 1756.9286
$ws.Range("J126").Value = 2185
$ws.Range("K126").Value = 5270.7858
$ws.Range("L126").Value = 6555
$ws.Range("M126").Value = -2800.7858
$ws.Range("N126").Value = -11495
$ws.Range("H132").Value = 10480842
$ws.Range("I132").Value = 14846946
$ws.Range("J132").Value = 2193.2666
$ws.Range("K132").Value = 44540838
$ws.Range("L132").Value = 6579.7998
$ws.Range("M132").Value = -44538308
$ws.Range("N132").Value = -11639.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 13200.333
$ws.Range("I14").Value = 12727.272
$ws.Range("J14").Value = 13474.211
$ws.Range("K14").Value = 12727.272
$ws.Range("L14").Value = 13474.211
$ws.Range("M14").Value = -12559.272
$ws.Range("N14").Value = -13810.211
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250
$ws.Range("H132").Value = 1192.909
$ws.Range("I132").Value = 869.08
$ws.Range("K132").Value = 2607.24
$ws.Range("M132").Value = -77.24000000000024
$ws.Range("H136").Value = 8931120
$ws.Range("I136").Value = 2677.35
$ws.Range("J136").Value = 31252224
$ws.Range("K136").Value = 8032.049999999999
$ws.Range("L136").Value = 93756672
$ws.Range("M136").Value = -5482.049999999999
$ws.Range("N136").Value = -93761772
